$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 23 new validation rows (new PB09 base station records) to the dataset.
$ws.Cells.Item(104, 1).Value = 'NuevaBase_459_PB09_TE'
$ws.Cells.Item(104, 2).Value = 'CX'
$ws.Cells.Item(104, 3).Value = 'PB09'
$ws.Cells.Item(104, 4).Value = 'HHE'
$ws.Cells.Item(104, 5).Value = '2019-03-28T15:04:21.000000Z'
$ws.Cells.Item(104, 6).Value = '2019-03-28T15:44:15.000000Z'
$ws.Cells.Item(105, 1).Value = 'NuevaBase_460_PB09_TE'
$ws.Cells.Item(105, 2).Value = 'CX'
$ws.Cells.Item(105, 3).Value = 'PB09'
$ws.Cells.Item(105, 4).Value = 'HHE'
$ws.Cells.Item(105, 5).Value = '2019-03-28T20:59:01.000001Z'
$ws.Cells.Item(105, 6).Value = '2019-03-28T21:08:01.000001Z'
$ws.Cells.Item(106, 1).Value = 'NuevaBase_462_PB09'
$ws.Cells.Item(106, 2).Value = 'CX'
$ws.Cells.Item(106, 3).Value = 'PB09'
$ws.Cells.Item(106, 4).Value = 'HHE'
$ws.Cells.Item(106, 5).Value = '2019-04-26T06:20:34.000000Z'
$ws.Cells.Item(106, 6).Value = '2019-04-26T06:28:34.000000Z'
$ws.Cells.Item(107, 1).Value = 'NuevaBase_463_PB09'
$ws.Cells.Item(107, 2).Value = 'CX'
$ws.Cells.Item(107, 3).Value = 'PB09'
$ws.Cells.Item(107, 4).Value = 'HHE'
$ws.Cells.Item(107, 5).Value = '2019-06-03T08:40:40.000000Z'
$ws.Cells.Item(107, 6).Value = '2019-06-03T08:45:52.000000Z'
$ws.Cells.Item(108, 1).Value = 'NuevaBase_464_PB09_TE'
$ws.Cells.Item(108, 2).Value = 'CX'
$ws.Cells.Item(108, 3).Value = 'PB09'
$ws.Cells.Item(108, 4).Value = 'HHE'
$ws.Cells.Item(108, 5).Value = '2019-06-14T00:17:32.000000Z'
$ws.Cells.Item(108, 6).Value = '2019-06-14T00:48:32.000000Z'
$ws.Cells.Item(109, 1).Value = 'NuevaBase_465_PB09_TE'
$ws.Cells.Item(109, 2).Value = 'CX'
$ws.Cells.Item(109, 3).Value = 'PB09'
$ws.Cells.Item(109, 4).Value = 'HHE'
$ws.Cells.Item(109, 5).Value = '2019-08-02T00:07:35.000000Z'
$ws.Cells.Item(109, 6).Value = '2019-08-02T00:38:35.000000Z'
$ws.Cells.Item(110, 1).Value = 'NuevaBase_467_PB09_TE'
$ws.Cells.Item(110, 2).Value = 'CX'
$ws.Cells.Item(110, 3).Value = 'PB09'
$ws.Cells.Item(110, 4).Value = 'HHE'
$ws.Cells.Item(110, 5).Value = '2019-09-26T20:31:08.000000Z'
$ws.Cells.Item(110, 6).Value = '2019-09-26T20:51:08.000000Z'
$ws.Cells.Item(111, 1).Value = 'NuevaBase_469_PB09'
$ws.Cells.Item(111, 2).Value = 'CX'
$ws.Cells.Item(111, 3).Value = 'PB09'
$ws.Cells.Item(111, 4).Value = 'HHE'
$ws.Cells.Item(111, 5).Value = '2019-11-04T21:52:45.000000Z'
$ws.Cells.Item(111, 6).Value = '2019-11-04T22:02:45.000000Z'
$ws.Cells.Item(112, 1).Value = 'NuevaBase_471_PB09_TE'
$ws.Cells.Item(112, 2).Value = 'CX'
$ws.Cells.Item(112, 3).Value = 'PB09'
$ws.Cells.Item(112, 4).Value = 'HHE'
$ws.Cells.Item(112, 5).Value = '2019-12-03T07:27:30.000000Z'
$ws.Cells.Item(112, 6).Value = '2019-12-03T07:44:00.000000Z'
$ws.Cells.Item(113, 1).Value = 'NuevaBase_472_PB09_TE'
$ws.Cells.Item(113, 2).Value = 'CX'
$ws.Cells.Item(113, 3).Value = 'PB09'
$ws.Cells.Item(113, 4).Value = 'HHE'
$ws.Cells.Item(113, 5).Value = '2019-12-03T08:45:46.000000Z'
$ws.Cells.Item(113, 6).Value = '2019-12-03T09:27:46.000000Z'
$ws.Cells.Item(114, 1).Value = 'NuevaBase_473_PB09_TE'
$ws.Cells.Item(114, 2).Value = 'CX'
$ws.Cells.Item(114, 3).Value = 'PB09'
$ws.Cells.Item(114, 4).Value = 'HHE'
$ws.Cells.Item(114, 5).Value = '2020-01-09T16:24:40.000000Z'
$ws.Cells.Item(114, 6).Value = '2020-01-09T16:44:40.000000Z'
$ws.Cells.Item(115, 1).Value = 'NuevaBase_475_PB09'
$ws.Cells.Item(115, 2).Value = 'CX'
$ws.Cells.Item(115, 3).Value = 'PB09'
$ws.Cells.Item(115, 4).Value = 'HHE'
$ws.Cells.Item(115, 5).Value = '2020-02-11T13:49:12.000000Z'
$ws.Cells.Item(115, 6).Value = '2020-02-11T13:58:30.000000Z'
$ws.Cells.Item(116, 1).Value = 'NuevaBase_476_PB09_TE'
$ws.Cells.Item(116, 2).Value = 'CX'
$ws.Cells.Item(116, 3).Value = 'PB09'
$ws.Cells.Item(116, 4).Value = 'HHE'
$ws.Cells.Item(116, 5).Value = '2020-02-13T05:49:38.000000Z'
$ws.Cells.Item(116, 6).Value = '2020-02-13T06:23:44.000000Z'
$ws.Cells.Item(117, 1).Value = 'NuevaBase_478_PB09_TE'
$ws.Cells.Item(117, 2).Value = 'CX'
$ws.Cells.Item(117, 3).Value = 'PB09'
$ws.Cells.Item(117, 4).Value = 'HHE'
$ws.Cells.Item(117, 5).Value = '2020-03-30T14:52:34.000000Z'
$ws.Cells.Item(117, 6).Value = '2020-03-30T15:06:16.000000Z'
$ws.Cells.Item(118, 1).Value = 'NuevaBase_479_PB09'
$ws.Cells.Item(118, 2).Value = 'CX'
$ws.Cells.Item(118, 3).Value = 'PB09'
$ws.Cells.Item(118, 4).Value = 'HHE'
$ws.Cells.Item(118, 5).Value = '2020-04-03T08:46:00.000000Z'
$ws.Cells.Item(118, 6).Value = '2020-04-03T08:55:54.000000Z'
$ws.Cells.Item(119, 1).Value = 'NuevaBase_480_PB09_TE'
$ws.Cells.Item(119, 2).Value = 'CX'
$ws.Cells.Item(119, 3).Value = 'PB09'
$ws.Cells.Item(119, 4).Value = 'HHE'
$ws.Cells.Item(119, 5).Value = '2020-06-15T04:33:44.000000Z'
$ws.Cells.Item(119, 6).Value = '2020-06-15T04:53:44.000000Z'
$ws.Cells.Item(120, 1).Value = 'NuevaBase_481_PB09_TE'
$ws.Cells.Item(120, 2).Value = 'CX'
$ws.Cells.Item(120, 3).Value = 'PB09'
$ws.Cells.Item(120, 4).Value = 'HHE'
$ws.Cells.Item(120, 5).Value = '2020-06-19T05:38:36.000000Z'
$ws.Cells.Item(120, 6).Value = '2020-06-19T06:02:24.000000Z'
$ws.Cells.Item(121, 1).Value = 'NuevaBase_493_PB09_TE'
$ws.Cells.Item(121, 2).Value = 'CX'
$ws.Cells.Item(121, 3).Value = 'PB09'
$ws.Cells.Item(121, 4).Value = 'HHE'
$ws.Cells.Item(121, 5).Value = '2020-10-06T05:15:39.000000Z'
$ws.Cells.Item(121, 6).Value = '2020-10-06T05:29:39.000000Z'
$ws.Cells.Item(122, 1).Value = 'NuevaBase_494_PB09_TE'
$ws.Cells.Item(122, 2).Value = 'CX'
$ws.Cells.Item(122, 3).Value = 'PB09'
$ws.Cells.Item(122, 4).Value = 'HHE'
$ws.Cells.Item(122, 5).Value = '2020-10-06T11:31:18.000000Z'
$ws.Cells.Item(122, 6).Value = '2020-10-06T11:50:18.000000Z'
$ws.Cells.Item(123, 1).Value = 'NuevaBase_496_PB09_TE'
$ws.Cells.Item(123, 2).Value = 'CX'
$ws.Cells.Item(123, 3).Value = 'PB09'
$ws.Cells.Item(123, 4).Value = 'HHE'
$ws.Cells.Item(123, 5).Value = '2020-10-21T09:11:36.000001Z'
$ws.Cells.Item(123, 6).Value = '2020-10-21T09:18:06.000001Z'
$ws.Cells.Item(124, 1).Value = 'NuevaBase_497_PB09_TE'
$ws.Cells.Item(124, 2).Value = 'CX'
$ws.Cells.Item(124, 3).Value = 'PB09'
$ws.Cells.Item(124, 4).Value = 'HHE'
$ws.Cells.Item(124, 5).Value = '2020-10-25T21:15:54.000000Z'
$ws.Cells.Item(124, 6).Value = '2020-10-25T21:42:30.000000Z'
$ws.Cells.Item(125, 1).Value = 'NuevaBase_498_PB09_TE'
$ws.Cells.Item(125, 2).Value = 'CX'
$ws.Cells.Item(125, 3).Value = 'PB09'
$ws.Cells.Item(125, 4).Value = 'HHE'
$ws.Cells.Item(125, 5).Value = '2020-10-28T14:51:40.000001Z'
$ws.Cells.Item(125, 6).Value = '2020-10-28T15:14:40.000001Z'
$ws.Cells.Item(126, 1).Value = 'NuevaBase_501_PB09'
$ws.Cells.Item(126, 2).Value = 'CX'
$ws.Cells.Item(126, 3).Value = 'PB09'
$ws.Cells.Item(126, 4).Value = 'HHE'
$ws.Cells.Item(126, 5).Value = '2020-11-21T02:10:10.000000Z'
$ws.Cells.Item(126, 6).Value = '2020-11-21T02:19:46.000000Z'

# Set column A width to fit the new (longer) station identifier strings.
$ws.Columns.Item(1).ColumnWidth = 21.6

# Move the active selection to cell A3 (matches the post-edit view state).
$ws.Range("A3").Select()
